$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 612, shifting the existing rows 612-655 down to 614-657.
$ws.Range("A612:R613").EntireRow.Insert()

# New row 612 data (based on old row 612, with D/J/K/L/M/O/P updated)
$ws.Cells.Item(612, 1).Value = 9
$ws.Cells.Item(612, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(612, 3).Value = "Metropolitana"
$ws.Cells.Item(612, 4).Value = 44931
$ws.Cells.Item(612, 5).Value = 13
$ws.Cells.Item(612, 6).Value = 100112031
$ws.Cells.Item(612, 7).Value = "Poroto verde"
$ws.Cells.Item(612, 8).Value = "Magnum"
$ws.Cells.Item(612, 9).Value = "Primera"
$ws.Cells.Item(612, 10).Value = 90
$ws.Cells.Item(612, 11).Value = 27000
$ws.Cells.Item(612, 12).Value = 30000
$ws.Cells.Item(612, 13).Value = 28500
$ws.Cells.Item(612, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(612, 15).Value = "Región Metropolitana"
$ws.Cells.Item(612, 16).Value = 1140
$ws.Cells.Item(612, 17).Value = 25
$ws.Cells.Item(612, 18).Value = "Hortaliza"

# New row 613 data (based on old row 613, with D/J/K/L/M/O/P updated)
$ws.Cells.Item(613, 1).Value = 9
$ws.Cells.Item(613, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(613, 3).Value = "Metropolitana"
$ws.Cells.Item(613, 4).Value = 44931
$ws.Cells.Item(613, 5).Value = 13
$ws.Cells.Item(613, 6).Value = 100112031
$ws.Cells.Item(613, 7).Value = "Poroto verde"
$ws.Cells.Item(613, 8).Value = "Sin especificar"
$ws.Cells.Item(613, 9).Value = "Primera"
$ws.Cells.Item(613, 10).Value = 52
$ws.Cells.Item(613, 11).Value = 30000
$ws.Cells.Item(613, 12).Value = 32000
$ws.Cells.Item(613, 13).Value = 31000
$ws.Cells.Item(613, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(613, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(613, 16).Value = 1240
$ws.Cells.Item(613, 17).Value = 25
$ws.Cells.Item(613, 18).Value = "Hortaliza"
